$d = $word.ActiveDocument

# 1) Remove the old "_GoBack" bookmark that was sitting at the very start of
#    the title paragraph. Saving renumbers the remaining bookmark ids
#    automatically (they shift down by one).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Expand "format and are 3.75 in x 2.75" into the fuller guidance about
#    using the U.S. Web Coated (SWOP) v2 color profile, and drop a fresh
#    "_GoBack" bookmark right after the newly inserted sentence (that is
#    where the cursor was left when the document was last saved).
$rng = $d.Content
$found = $rng.Find.Execute("format and are 3.75 in x 2.75")
if ($found) {
    $segStart = $rng.Start
    $newText = "format, use U.S. Web Coated (SWOP) v2 profile, and have 3.75 in x 2.75"
    $target = $d.Range($segStart, $rng.End)
    $target.Text = $newText

    $markerOffset = "format, use U.S. Web Coated (SWOP) v2 profile,".Length
    $markerPos = $segStart + $markerOffset
    $bmRange = $d.Range($markerPos, $markerPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# 3) Tidy up two spots where an adjoining run split was left over from an
#    earlier edit: merge the pair of runs back into a single run. Re-assign
#    the exact existing text (round-tripped through a throwaway edit so the
#    engine actually re-serializes the span) so no characters change.
function Merge-Runs($rangeStart, $rangeEnd) {
    $mrng = $d.Range($rangeStart, $rangeEnd)
    $original = $mrng.Text
    $mrng.Text = "."
    $tmp = $d.Range($rangeStart, $rangeStart + 1)
    $tmp.Text = $original
}

$rng2 = $d.Content
if ($rng2.Find.Execute("- A customized pack")) {
    Merge-Runs $rng2.Start $rng2.End
}

$rng3 = $d.Content
if ($rng3.Find.Execute("You can find the PDF at: ")) {
    # only merge the "the" + " PDF at: " pair, not the whole sentence
    $mergeEnd = $rng3.End
    $mergeStart = $mergeEnd - "the PDF at: ".Length
    Merge-Runs $mergeStart $mergeEnd
}
